$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.228358268737793
$ws.Range("B1").Value = 2.433889389038086
$ws.Range("C1").Value = 2.344024896621704
$ws.Range("D1").Value = 2.666018486022949
$ws.Range("E1").Value = 0.457331120967865
